# Auto-generated edit script applying the Twintania_Profits leve-profit recompute.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the
# affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1858.5264
$ws.Range("J40").Value = 1933.3334
$ws.Range("L40").Value = 1933.3334
$ws.Range("N40").Value = -2283.3334
$ws.Range("H51").Value = 7999.1665
$ws.Range("I51").Value = 5999.3335
$ws.Range("K51").Value = 5999.3335
$ws.Range("M51").Value = -5515.3335
$ws.Range("H61").Value = 6789.25
$ws.Range("I61").Value = 862.8
$ws.Range("J61").Value = 16666.666
$ws.Range("K61").Value = 2588.4
$ws.Range("L61").Value = 49999.99800000001
$ws.Range("M61").Value = -2416.4
$ws.Range("N61").Value = -50343.99800000001
$ws.Range("H98").Value = 1477.7142
$ws.Range("I98").Value = 932.3333
$ws.Range("K98").Value = 932.3333
$ws.Range("M98").Value = 565.6667
$ws.Range("H106").Value = 26752.055
$ws.Range("I106").Value = 11579.143
$ws.Range("J106").Value = 36407.547
$ws.Range("K106").Value = 11579.143
$ws.Range("L106").Value = 36407.547
$ws.Range("M106").Value = -10948.143
$ws.Range("N106").Value = -37669.547
$ws.Range("H112").Value = 2182.6667
$ws.Range("I112").Value = 1883.5
$ws.Range("J112").Value = 2201.9678
$ws.Range("K112").Value = 5650.5
$ws.Range("L112").Value = 6605.903399999999
$ws.Range("M112").Value = -4542.5
$ws.Range("N112").Value = -8821.903399999999
$ws.Range("H116").Value = 6891.6665
$ws.Range("I116").Value = 4668.75
$ws.Range("K116").Value = 4668.75
$ws.Range("M116").Value = -1226.75
$ws.Range("H122").Value = 1477.7142
$ws.Range("I122").Value = 932.3333
$ws.Range("K122").Value = 2796.9999
$ws.Range("M122").Value = -346.9998999999998
$ws.Range("H132").Value = 1756.5807
$ws.Range("I132").Value = 1479.1538
$ws.Range("K132").Value = 4437.4614
$ws.Range("M132").Value = -1907.4614
$ws.Range("H135").Value = 2082.2068
$ws.Range("I135").Value = 1981.5714
$ws.Range("J135").Value = 4900
$ws.Range("K135").Value = 17834.1426
$ws.Range("L135").Value = 44100
$ws.Range("M135").Value = -15299.1426
$ws.Range("N135").Value = -49170

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3618.8572
$ws.Range("I32").Value = 2947.158
$ws.Range("K32").Value = 2947.158
$ws.Range("M32").Value = -2660.158
$ws.Range("H44").Value = 15969
$ws.Range("J44").Value = 15969
$ws.Range("L44").Value = 15969
$ws.Range("N44").Value = -16945
$ws.Range("H55").Value = 23676.938
$ws.Range("J55").Value = 29530.455
$ws.Range("L55").Value = 29530.455
$ws.Range("N55").Value = -30160.455
$ws.Range("H80").Value = 32861
$ws.Range("J80").Value = 32861
$ws.Range("L80").Value = 32861
$ws.Range("N80").Value = -34857
$ws.Range("H83").Value = 32861
$ws.Range("J83").Value = 32861
$ws.Range("L83").Value = 98583
$ws.Range("N83").Value = -108567
$ws.Range("H102").Value = 1339.8846
$ws.Range("I102").Value = 1339.8846
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1339.8846
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = 282.1153999999999
$ws.Range("H132").Value = 4668.727
$ws.Range("I132").Value = 4660.6924
$ws.Range("K132").Value = 13982.0772
$ws.Range("M132").Value = -11452.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 28036.79
$ws.Range("I82").Value = 9034.333000000001
$ws.Range("J82").Value = 45139
$ws.Range("K82").Value = 9034.333000000001
$ws.Range("L82").Value = 45139
$ws.Range("M82").Value = -8651.333000000001
$ws.Range("N82").Value = -45905
$ws.Range("H85").Value = 28036.79
$ws.Range("I85").Value = 9034.333000000001
$ws.Range("J85").Value = 45139
$ws.Range("K85").Value = 9034.333000000001
$ws.Range("L85").Value = 45139
$ws.Range("M85").Value = -7708.333000000001
$ws.Range("N85").Value = -47791

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2979.9565
$ws.Range("I31").Value = 1927.6666
$ws.Range("J31").Value = 3137.8
$ws.Range("K31").Value = 1927.6666
$ws.Range("L31").Value = 3137.8
$ws.Range("M31").Value = -1632.6666
$ws.Range("N31").Value = -3727.8
$ws.Range("H34").Value = 2979.9565
$ws.Range("I34").Value = 1927.6666
$ws.Range("J34").Value = 3137.8
$ws.Range("K34").Value = 1927.6666
$ws.Range("L34").Value = 3137.8
$ws.Range("M34").Value = -1725.6666
$ws.Range("N34").Value = -3541.8
$ws.Range("H86").Value = 3044.4614
$ws.Range("J86").Value = 3083
$ws.Range("L86").Value = 3083
$ws.Range("N86").Value = -5329
$ws.Range("H89").Value = 3044.4614
$ws.Range("J89").Value = 3083
$ws.Range("L89").Value = 15415
$ws.Range("N89").Value = -26647
$ws.Range("H132").Value = 28099.47
$ws.Range("I132").Value = 17174.734
$ws.Range("K132").Value = 51524.202
$ws.Range("M132").Value = -48994.202
$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1550.5264
$ws.Range("I5").Value = 1075.9231
$ws.Range("J5").Value = 1797.32
$ws.Range("K5").Value = 3227.7693
$ws.Range("L5").Value = 5391.96
$ws.Range("M5").Value = -3115.7693
$ws.Range("N5").Value = -5615.96
$ws.Range("H114").Value = 2420
$ws.Range("J114").Value = 2923.5
$ws.Range("L114").Value = 8770.5
$ws.Range("N114").Value = -15278.5
$ws.Range("H132").Value = 28573212
$ws.Range("J132").Value = 33334862
$ws.Range("L132").Value = 300013758
$ws.Range("N132").Value = -300018818
$ws.Range("H135").Value = 1550.5264
$ws.Range("I135").Value = 1075.9231
$ws.Range("J135").Value = 1797.32
$ws.Range("K135").Value = 9683.3079
$ws.Range("L135").Value = 16175.88
$ws.Range("M135").Value = -7148.3079
$ws.Range("N135").Value = -21245.88
$ws.Range("H139").Value = 1372.2858
$ws.Range("I139").Value = 1372.2858
$ws.Range("K139").Value = 4116.857400000001
$ws.Range("M139").Value = 1023.142599999999
$ws.Range("H140").Value = 5042.615
$ws.Range("I140").Value = 4713.25
$ws.Range("J140").Value = 8995
$ws.Range("K140").Value = 14139.75
$ws.Range("L140").Value = 26985
$ws.Range("M140").Value = -8959.75
$ws.Range("N140").Value = -37345

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9147
$ws.Range("I70").Value = 8135.625
$ws.Range("K70").Value = 8135.625
$ws.Range("M70").Value = -7865.625
$ws.Range("H73").Value = 9147
$ws.Range("I73").Value = 8135.625
$ws.Range("K73").Value = 8135.625
$ws.Range("M73").Value = -7199.625
$ws.Range("H97").Value = 2447.3
$ws.Range("I97").Value = 2053.1428
$ws.Range("K97").Value = 2053.1428
$ws.Range("M97").Value = -1557.1428
$ws.Range("H102").Value = 6506.5
$ws.Range("I102").Value = 6506.5
$ws.Range("K102").Value = 6506.5
$ws.Range("M102").Value = -4884.5
$ws.Range("H122").Value = 3056
$ws.Range("I122").Value = 2967.2
$ws.Range("K122").Value = 8901.599999999999
$ws.Range("M122").Value = -6451.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 719.5
$ws.Range("I22").Value = 622.875
$ws.Range("K22").Value = 622.875
$ws.Range("M22").Value = -327.875
$ws.Range("H27").Value = 719.5
$ws.Range("I27").Value = 622.875
$ws.Range("K27").Value = 622.875
$ws.Range("M27").Value = -515.875
$ws.Range("H46").Value = 1173.4131
$ws.Range("I46").Value = 807.6875
$ws.Range("K46").Value = 807.6875
$ws.Range("M46").Value = -619.6875
$ws.Range("H93").Value = 4492.875
$ws.Range("I93").Value = 4449.067
$ws.Range("K93").Value = 4449.067
$ws.Range("M93").Value = -3201.067
$ws.Range("H136").Value = 5351.763
$ws.Range("I136").Value = 5099.8184
$ws.Range("J136").Value = 7014.6
$ws.Range("K136").Value = 15299.4552
$ws.Range("L136").Value = 21043.8
$ws.Range("M136").Value = -12749.4552
$ws.Range("N136").Value = -26143.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 10281.667
$ws.Range("I96").Value = 12591.4
$ws.Range("J96").Value = 7394.5
$ws.Range("K96").Value = 12591.4
$ws.Range("L96").Value = 7394.5
$ws.Range("M96").Value = -11218.4
$ws.Range("N96").Value = -10140.5
